$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab
$ws.Name = "Sprint 3"

# Add new row content
$ws.Range("A2").Value = "Team should begin working on Final Presentation for the Final Project deliverable"
$ws.Range("B2").Value = "N/A"
$ws.Range("C2").Value = "N.A"

# Style the new row similarly to the header data rows (wrap text)
$ws.Range("A2:C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 30

# A2 gets its own font style: Arial 12 black, wrap text only (no vertical centering)
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 12
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").VerticalAlignment = -4107  # xlBottom (clears the inherited vertical=center)

# B2/C2 reuse the existing column style, same Arial 12 vertical-centered wrap
$ws.Range("B2:C2").Font.Name = "Arial"
$ws.Range("B2:C2").Font.Size = 12
$ws.Range("B2:C2").VerticalAlignment = -4108  # xlCenter

# Update selection to A3 as shown in diff
$ws.Range("A3").Select()
